# fx_predict.xlsx update ("upd fx_predict.xlsx by request")
# Updates the FX forecast table on sheet "Предсказание" with refreshed figures
# and period headers, then leaves the selection on G8 as in the authored file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): quarter labels shift forward one period ---
# D1: "3К23" -> "4К23", E1: 2023 (number) -> "1к24", F1: "1к24" -> "2к24"
$ws.Range("A2").Value = "USD/RUB"   # drop leading space, set before F1/D1 so new
$ws.Range("E1").Value = "1к24"      # shared-string entries are appended in the
$ws.Range("F1").Value = "2к24"      # same order as the authored workbook
$ws.Range("D1").Value = "4К23"

# --- Row 2: USD/RUB ---
$ws.Range("B2").Value = 98
$ws.Range("C2").Value = 95
$ws.Range("E2").Value = 90
$ws.Range("F2").Value = 90

# --- Row 3: EUR/USD ---
$ws.Range("B3").Value = 1.07
$ws.Range("C3").Value = 1.07
$ws.Range("D3").Value = 1.07
$ws.Range("E3").Value = 1.07
$ws.Range("F3").Value = 1.078

# --- Row 4: EUR/RUB ---
$ws.Range("B4").Value = 104.86
$ws.Range("C4").Value = 101.65
$ws.Range("D4").Value = 96.3
$ws.Range("E4").Value = 96.3
$ws.Range("F4").Value = 97.2

# --- Row 5: USD/CNY ---
$ws.Range("B5").Value = 7.3
$ws.Range("C5").Value = 7.3
$ws.Range("E5").Value = 7.1
$ws.Range("F5").Value = 7

# --- Row 6: CNY/RUB ---
$ws.Range("B6").Value = 13.42
$ws.Range("C6").Value = 13.01
$ws.Range("E6").Value = 12.68
$ws.Range("F6").Value = 12.86

# --- Row 7: USD/INR ---
$ws.Range("B7").Value = 83
$ws.Range("C7").Value = 83
$ws.Range("D7").Value = 83
$ws.Range("E7").Value = 83
$ws.Range("F7").Value = 83

# --- Row 8: INR/RUB ---
$ws.Range("B8").Value = 1.18
$ws.Range("C8").Value = 1.14
$ws.Range("D8").Value = 1.08
$ws.Range("E8").Value = 1.08
$ws.Range("F8").Value = 1.08

# --- Row 9: USD/TRY ---
$ws.Range("B9").Value = 27.5
$ws.Range("C9").Value = 27.5
$ws.Range("D9").Value = 28
$ws.Range("E9").Value = 29
$ws.Range("F9").Value = 30

# --- Row 10: TRY/RUB ---
$ws.Range("B10").Value = 3.56
$ws.Range("C10").Value = 3.45
$ws.Range("D10").Value = 3.21
$ws.Range("E10").Value = 3.1
$ws.Range("F10").Value = 3

# --- Row 11: USD/KZT ---
$ws.Range("B11").Value = 480
$ws.Range("C11").Value = 480
$ws.Range("D11").Value = 480
$ws.Range("E11").Value = 480
$ws.Range("F11").Value = 480

# --- Row 12: KZT/RUB*100 ---
$ws.Range("B12").Value = 20.42
$ws.Range("C12").Value = 19.79
$ws.Range("D12").Value = 18.75
$ws.Range("E12").Value = 18.75
$ws.Range("F12").Value = 18.75

# Restore the authored selection (G8)
$ws.Range("G8").Select()
